# E-Learning Project.xlsx - add .NET Identity/JWT authentication task rows
# (commit: "Add user authentication with jwt. add controller action Register,
#  Login, RefreshToken and RevokeToken")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("B1").Value = "Acceptance Criteria"

# --- Row 2: Setup Web Api project using .Net 5.0 ------------------------
$ws.Range("A2").Value = "Setup Web Api project using .Net 5.0"
$ws.Range("B2").Value = "Install Microsoft.EntityFrameworkCore`nInstall Microsoft.EntityFrameworkCore.Design`nInstall Microsoft.EntityFrameworkCore.Tools`nInstall Npgsql.EntityFrameworkCore.PostgreSql`nInstall Npgsql`nInstall Microsoft.AspNetCore.Identity.EntityFrameworkCore"
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = "Finished"
$ws.Rows.Item(2).RowHeight = 86.4

# --- Row 3: Setup Configuration ------------------------------------------
$ws.Range("A3").Value = "Setup Configuration"
$ws.Range("B3").Value = "Add SqlConnection to appSetting.json, Add IdentityDbContext, Setup Configuration at Startup.cs"
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = "Finished"
$ws.Rows.Item(3).RowHeight = 28.8

# --- Row 4: Setup Swagger and configure ----------------------------------
$ws.Range("A4").Value = "Setup Swagger and configure"
$ws.Range("B4").Value = "Add Swashbuckle.AspNetCore.Swagger, Swashbuckle.AspNetCore.SwaggerGen, Swashbuckle.AspNetCore.SwaggerUI"
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "Finished"
$ws.Rows.Item(4).RowHeight = 43.2

# --- Row 5: Setup xml Documentation ---------------------------------------
$ws.Range("A5").Value = "Setup xml Documentation"
$ws.Range("B5").Value = "Add xml Documentation at properties -> build of the project. Supress warning and errors if need"
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "Finished"
$ws.Rows.Item(5).RowHeight = 28.8

# --- Row 6: Setup Logging --------------------------------------------------
$ws.Range("A6").Value = "Setup Logging"
$ws.Range("B6").Value = "add NLog.Extensions.Logging,`nadd nlog.config file at root,`nadd LoggerService and configure it at Startup.cs"
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = "Finished"
$ws.Rows.Item(6).RowHeight = 43.2

# --- Row 7: Setup Cors Policy -----------------------------------------------
$ws.Range("A7").Value = "Setup Cors Policy"
$ws.Range("B7").Value = "Add Cors Policy to interact with API with frontend"
$ws.Range("C7").Value = $null
$ws.Range("D7").Value = "Finished"
$ws.Rows.Item(7).RowHeight = 14.4

# --- Row 8: Setup Git ---------------------------------------------------------
$ws.Range("A8").Value = "Setup Git"
$ws.Range("B8").Value = "Use Git and Github"
$ws.Range("C8").Value = $null
$ws.Range("D8").Value = "Finished"
$ws.Rows.Item(8).RowHeight = 14.4

# --- Row 9: Setup Identity Schema and Migrate (NEW) ---------------------------
$ws.Range("A9").Value = "Setup Identity Schema and Migrate`nwith EntityFrameworkCore"
$ws.Range("B9").Value = "Add AppUser<int> and RoleUser<int> classes and update DbContext to`nIdentityDbContext<AppUser, AppRole, int>"
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = "Finished"
$ws.Rows.Item(9).RowHeight = 43.2

# --- Row 10: Configure Authentication and Identity (NEW) ---------------------
$ws.Range("A10").Value = "Configure Authentication and Identity`nand authentication with JWT"
$ws.Range("B10").Value = "Seed the sample User Data, Create UserController,`nAdd controller actions for Signup, Login, Refrsh Token and Revoke Token`nGenerate JWT, Set key and Issuer at appSetting file,`nInstall Microsoft.AspNetCore.Authentication.JwtBearer"
$ws.Range("C10").Value = "Working"
$ws.Range("D10").Value = $null
$ws.Rows.Item(10).RowHeight = 72

# Apply the alignment/wrap formatting that the A/B columns use throughout
# (vertical-top, wrapped text) to the two newly-added rows.
$ws.Range("A9:B10").VerticalAlignment = -4160
$ws.Range("A9:B10").WrapText = $true
$ws.Range("C9:D10").VerticalAlignment = -4160
$ws.Range("C9:D10").WrapText = $false

# --- New trailing blank rows 14 & 15 (sheet grew from 13 to 15 rows) ---------
$ws.Range("A14:D15").VerticalAlignment = -4160
$ws.Range("A14:D15").WrapText = $false

# --- Selection / view state to match the saved workbook ----------------------
$ws.Range("B11").Select()
